$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the value/percentage columns so numeric-looking
# strings (e.g. "219.40", "0.0847") are not reinterpreted as numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "26.871.19"
$ws.Range("E2").Value = "  +0.49%  "
$ws.Range("D3").Value = "1.641.77"
$ws.Range("E3").Value = "  -0.10%  "
$ws.Range("E4").Value = "  -0.44%  "
$ws.Range("D5").Value = "219.40"
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  -0.22%  "
$ws.Range("E8").Value = "  -0.18%  "
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("D10").Value = "19.26"
$ws.Range("E10").Value = "  +0.82%  "
$ws.Range("D11").Value = "0.0847"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("D12").Value = "1.871.07"
$ws.Range("E12").Value = "  +0.06%  "
$ws.Range("D13").Value = "1.641.65"
$ws.Range("E13").Value = "  +0.07%  "
$ws.Range("E14").Value = "  -0.54%  "
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").Value = "65.43"
$ws.Range("E16").Value = "  +1.30%  "
$ws.Range("D17").Value = "26.869.78"
$ws.Range("E17").Value = "  +0.54%  "
$ws.Range("E18").Value = "  -0.42%  "
$ws.Range("D19").Value = "216.59"
$ws.Range("E19").Value = "  +1.19%  "
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("E21").Value = "  -0.08%  "
$ws.Range("D22").Value = "6.57"
$ws.Range("E22").Value = "  +5.12%  "
$ws.Range("D23").Value = "2.41"
$ws.Range("E23").Value = "  -2.01%  "
$ws.Range("D24").Value = "9.20"
$ws.Range("E24").Value = "  -1.62%  "
$ws.Range("D25").Value = "147.91"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +1.47%  "
$ws.Range("D29").Value = "15.78"
$ws.Range("E29").Value = "  +0.83%  "
$ws.Range("D30").Value = "0.0510"
$ws.Range("E30").Value = "  +0.36%  "
$ws.Range("D31").Value = "1.19"
$ws.Range("E31").Value = "  +1.01%  "
$ws.Range("E32").Value = "  -0.24%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("E34").Value = "  +1.71%  "
$ws.Range("D35").Value = "1.271.82"
$ws.Range("E35").Value = "  -1.47%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  -1.99%  "
$ws.Range("D38").Value = "0.531"
$ws.Range("E38").Value = "  -0.66%  "
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("E40").Value = "  -0.20%  "
$ws.Range("D41").Value = "0.805"
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").Value = "1.781.71"
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("D44").Value = "92.61"
$ws.Range("E44").Value = "  +1.32%  "
$ws.Range("D45").Value = "60.91"
$ws.Range("E45").Value = "  -0.78%  "
$ws.Range("D46").Value = "2.07"
$ws.Range("E46").Value = "  -7.60%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0515"
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "7.60"
$ws.Range("E49").Value = "  -1.20%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "0.0970"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").Value = "1.01"
$ws.Range("E51").Value = "  -0.13%  "

# Restore the original (default) style so no residual formatting is left behind.
$ws.Range("B2:E51").Style = "Normal"
